$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect it so the cell values can be updated,
# then re-apply protection afterwards so the sheet remains protected.
$ws.Unprotect()

# Update the confidential disclaimer date (A10) from 2021-03-25 to 2021-03-26
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2453512078365138
$ws.Range("E2").Value = 0.01918294849023083

$ws.Range("D3").Value = 0.4990503163980574
$ws.Range("E3").Value = 0.01843060597598423

$ws.Range("D4").Value = 0.09775358029011641
$ws.Range("E4").Value = 0.0222772277227723

$ws.Range("D5").Value = 0.1002972624064649
$ws.Range("E5").Value = 0.01112737920937046

$ws.Range("D6").Value = 0.05754763306884758
$ws.Range("E6").Value = 0.0186584789796882

# Update Percent Change (E) for Total row (row 7); D7 (weight) remains 1 (unchanged)
$ws.Range("E7").Value = 0.01827183506891017

# Restore sheet protection
$ws.Protect()
